# Insert 3 new data rows into the "Hortaliza, Vega Central Mapocho de Santiago - Zanahoria"
# sheet, right before the existing row 1327. This pushes the existing rows 1327-1393 down
# to 1330-1396 (dimension grows from A1:R1393 to A1:R1396), matching three new weekly
# price observations that were added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1327:1396 down by inserting 3 blank rows at 1327.
$ws.Range("1327:1329").Insert()

# Common/shared values for this data set (identical across every row in the sheet).
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$categoriaId = 100114013
$categoria = "Zanahoria"
$variedad  = "Sin especificar"
$unidad    = "$/saco 20 kilos"
$kgUnidades = 20
$clasificacion = "Hortaliza"
$fecha = 44939

# New row 1327: Camote
$ws.Cells.Item(1327, 1).Value = $mercadoId
$ws.Cells.Item(1327, 2).Value = $mercado
$ws.Cells.Item(1327, 3).Value = $region
$ws.Cells.Item(1327, 4).Value = $fecha
$ws.Cells.Item(1327, 5).Value = $codreg
$ws.Cells.Item(1327, 6).Value = $categoriaId
$ws.Cells.Item(1327, 7).Value = $categoria
$ws.Cells.Item(1327, 8).Value = $variedad
$ws.Cells.Item(1327, 9).Value = "Camote"
$ws.Cells.Item(1327, 10).Value = 250
$ws.Cells.Item(1327, 11).Value = 10000
$ws.Cells.Item(1327, 12).Value = 10000
$ws.Cells.Item(1327, 13).Value = 10000
$ws.Cells.Item(1327, 14).Value = $unidad
$ws.Cells.Item(1327, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1327, 16).Value = 500
$ws.Cells.Item(1327, 17).Value = $kgUnidades
$ws.Cells.Item(1327, 18).Value = $clasificacion

# New row 1328: Primera
$ws.Cells.Item(1328, 1).Value = $mercadoId
$ws.Cells.Item(1328, 2).Value = $mercado
$ws.Cells.Item(1328, 3).Value = $region
$ws.Cells.Item(1328, 4).Value = $fecha
$ws.Cells.Item(1328, 5).Value = $codreg
$ws.Cells.Item(1328, 6).Value = $categoriaId
$ws.Cells.Item(1328, 7).Value = $categoria
$ws.Cells.Item(1328, 8).Value = $variedad
$ws.Cells.Item(1328, 9).Value = "Primera"
$ws.Cells.Item(1328, 10).Value = 430
$ws.Cells.Item(1328, 11).Value = 12000
$ws.Cells.Item(1328, 12).Value = 12000
$ws.Cells.Item(1328, 13).Value = 12000
$ws.Cells.Item(1328, 14).Value = $unidad
$ws.Cells.Item(1328, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1328, 16).Value = 600
$ws.Cells.Item(1328, 17).Value = $kgUnidades
$ws.Cells.Item(1328, 18).Value = $clasificacion

# New row 1329: Segunda
$ws.Cells.Item(1329, 1).Value = $mercadoId
$ws.Cells.Item(1329, 2).Value = $mercado
$ws.Cells.Item(1329, 3).Value = $region
$ws.Cells.Item(1329, 4).Value = $fecha
$ws.Cells.Item(1329, 5).Value = $codreg
$ws.Cells.Item(1329, 6).Value = $categoriaId
$ws.Cells.Item(1329, 7).Value = $categoria
$ws.Cells.Item(1329, 8).Value = $variedad
$ws.Cells.Item(1329, 9).Value = "Segunda"
$ws.Cells.Item(1329, 10).Value = 340
$ws.Cells.Item(1329, 11).Value = 11000
$ws.Cells.Item(1329, 12).Value = 11000
$ws.Cells.Item(1329, 13).Value = 11000
$ws.Cells.Item(1329, 14).Value = $unidad
$ws.Cells.Item(1329, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1329, 16).Value = 550
$ws.Cells.Item(1329, 17).Value = $kgUnidades
$ws.Cells.Item(1329, 18).Value = $clasificacion

# Ensure the date column keeps the date/time number format used elsewhere in column D.
$ws.Range("D1327:D1329").NumberFormat = "YYYY-MM-DD HH:MM:SS"
